$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 10.06794800843628
$ws.Range("D2").Value = 3.983777410211226
$ws.Range("E2").Value = 12.16657276785682
$ws.Range("F2").Value = 23.38826680322079
$ws.Range("G2").Value = 28.01292897430985
$ws.Range("H2").Value = 13.49685549273043
$ws.Range("I2").Value = 21.67762619005112
$ws.Range("L2").Value = 9.176037176907197
$ws.Range("M2").Value = 21.43271464196771
$ws.Range("N2").Value = 17.44077157824414
$ws.Range("O2").Value = 20.55991972793914
$ws.Range("C3").Value = 10.11295619119016
$ws.Range("D3").Value = 3.968187591908877
$ws.Range("E3").Value = 12.2450925967679
$ws.Range("F3").Value = 23.19104253171543
$ws.Range("G3").Value = 27.5794044171809
$ws.Range("H3").Value = 13.49644824827044
$ws.Range("I3").Value = 21.63054172415454
$ws.Range("L3").Value = 9.221551313212375
$ws.Range("M3").Value = 20.74802168247708
$ws.Range("N3").Value = 17.14173113664868
$ws.Range("O3").Value = 20.47873960256266
$ws.Range("C4").Value = 10.14270205021584
$ws.Range("D4").Value = 3.958404878136829
$ws.Range("E4").Value = 12.29575262194584
$ws.Range("F4").Value = 23.07648756809865
$ws.Range("G4").Value = 27.32017458332991
$ws.Range("H4").Value = 13.49914079998883
$ws.Range("I4").Value = 21.60763397498322
$ws.Range("L4").Value = 9.250837214704529
$ws.Range("M4").Value = 20.31555881327071
$ws.Range("N4").Value = 16.957348357152
$ws.Range("O4").Value = 20.43453689771981
$ws.Range("C5").Value = 10.15535349392925
$ws.Range("D5").Value = 3.954366250305936
$ws.Range("E5").Value = 12.31701374010306
$ws.Range("F5").Value = 23.03149930035202
$ws.Range("G5").Value = 27.21643793671273
$ws.Range("H5").Value = 13.50097818266113
$ws.Range("I5").Value = 21.59981382905633
$ws.Range("L5").Value = 9.263109372502807
$ws.Range("M5").Value = 20.13652424627947
$ws.Range("N5").Value = 16.8821102927801
$ws.Range("O5").Value = 20.41795601852297
$ws.Range("C6").Value = 10.15748622472889
$ws.Range("D6").Value = 3.953692517515762
$ws.Range("E6").Value = 12.32058141356637
$ws.Range("F6").Value = 23.02413269203826
$ws.Range("G6").Value = 27.19933177333813
$ws.Range("H6").Value = 13.5013279798488
$ws.Range("I6").Value = 21.59860693080307
$ws.Range("L6").Value = 9.265167592175326
$ws.Range("M6").Value = 20.1066335503984
$ws.Range("N6").Value = 16.86961375502937
$ws.Range("O6").Value = 20.41528961927258
$ws.Range("C7").Value = 10.14287052829476
$ws.Range("D7").Value = 3.958350621610342
$ws.Range("E7").Value = 12.29603685825273
$ws.Range("F7").Value = 23.07587392056019
$ws.Range("G7").Value = 27.31876765667494
$ws.Range("H7").Value = 13.49916258280363
$ws.Range("I7").Value = 21.6075223696337
$ws.Range("L7").Value = 9.251001351858157
$ws.Range("M7").Value = 20.31315531879466
$ws.Range("N7").Value = 16.95633395303034
$ws.Range("O7").Value = 20.4343074677166
$ws.Range("C8").Value = 10.08302817535596
$ws.Range("D8").Value = 3.978446232162941
$ws.Range("E8").Value = 12.19313884053143
$ws.Range("F8").Value = 23.31893659797687
$ws.Range("G8").Value = 27.86209861418657
$ws.Range("H8").Value = 13.49610457475403
$ws.Range("I8").Value = 21.66014955120233
$ws.Range("L8").Value = 9.191452902203398
$ws.Range("M8").Value = 21.19926682027574
$ws.Range("N8").Value = 17.33787771052577
$ws.Range("O8").Value = 20.53076636527592
$ws.Range("C9").Value = 9.982468345645799
$ws.Range("D9").Value = 4.016147469150274
$ws.Range("E9").Value = 12.01073560617426
$ws.Range("F9").Value = 23.84516963075634
$ws.Range("G9").Value = 28.97590370799868
$ws.Range("H9").Value = 13.51343170223039
$ws.Range("I9").Value = 21.81065615412383
$ws.Range("L9").Value = 9.085266942081539
$ws.Range("M9").Value = 22.83218543524456
$ws.Range("N9").Value = 18.07598472000492
$ws.Range("O9").Value = 20.76401302003138
$ws.Range("C10").Value = 9.9188833007121
$ws.Range("D10").Value = 4.042759620274968
$ws.Range("E10").Value = 11.88846902032672
$ws.Range("F10").Value = 24.25871120860701
$ws.Range("G10").Value = 29.814401820609
$ws.Range("H10").Value = 13.54032536750759
$ws.Range("I10").Value = 21.94953594273148
$ws.Range("L10").Value = 9.013644026840749
$ws.Range("M10").Value = 23.9572964471125
$ws.Range("N10").Value = 18.60697677394123
$ws.Range("O10").Value = 20.96122143858041
$ws.Range("C11").Value = 9.892205776273968
$ws.Range("D11").Value = 4.054620291681417
$ws.Range("E11").Value = 11.835382770823
$ws.Range("F11").Value = 24.45193321472672
$ws.Range("G11").Value = 30.19828522921716
$ws.Range("H11").Value = 13.55561499718405
$ws.Range("I11").Value = 22.01871372130926
$ws.Range("L11").Value = 8.982435368420564
$ws.Range("M11").Value = 24.45104793376569
$ws.Range("N11").Value = 18.84506519655215
$ws.Range("O11").Value = 21.05628912067429
$ws.Range("C12").Value = 9.882428060583596
$ws.Range("D12").Value = 4.059075478921057
$ws.Range("E12").Value = 11.81564372998795
$ws.Range("F12").Value = 24.52576911941575
$ws.Range("G12").Value = 30.3438446198561
$ws.Range("H12").Value = 13.56184168015923
$ws.Range("I12").Value = 22.04575747226358
$ws.Range("L12").Value = 8.97081389515424
$ws.Range("M12").Value = 24.63528608183375
$ws.Range("N12").Value = 18.93464526715811
$ws.Range("O12").Value = 21.09303469009348
$ws.Range("C13").Value = 9.884519414617396
$ws.Range("D13").Value = 4.058117598763675
$ws.Range("E13").Value = 11.81987872707544
$ws.Range("F13").Value = 24.5098386426314
$ws.Range("G13").Value = 30.31248981226089
$ws.Range("H13").Value = 13.56048126964786
$ws.Range("I13").Value = 22.03989567695558
$ws.Range("L13").Value = 8.973308058584758
$ws.Range("M13").Value = 24.59573074140418
$ws.Range("N13").Value = 18.91537958357364
$ws.Range("O13").Value = 21.08508812804985
$ws.Range("C14").Value = 9.891394847215592
$ws.Range("D14").Value = 4.054987553168463
$ws.Range("E14").Value = 11.83375154566122
$ws.Range("F14").Value = 24.45799475689119
$ws.Range("G14").Value = 30.21025756267059
$ws.Range("H14").Value = 13.55611853155634
$ws.Range("I14").Value = 22.02092173196439
$ws.Range("L14").Value = 8.98147532821444
$ws.Range("M14").Value = 24.46626085395831
$ws.Range("N14").Value = 18.85244705218472
$ws.Range("O14").Value = 21.05929740699691
$ws.Range("C15").Value = 9.895648543718359
$ws.Range("D15").Value = 4.053065569606083
$ws.Range("E15").Value = 11.84229637390702
$ws.Range("F15").Value = 24.42632374443394
$ws.Range("G15").Value = 30.14765754945519
$ws.Range("H15").Value = 13.55350303367762
$ws.Range("I15").Value = 22.00940955052582
$ws.Range("L15").Value = 8.986503590210198
$ws.Range("M15").Value = 24.38659686758263
$ws.Range("N15").Value = 18.81382134753561
$ws.Range("O15").Value = 21.04359617130393
$ws.Range("C16").Value = 9.920672076627453
$ws.Range("D16").Value = 4.041979482662519
$ws.Range("E16").Value = 11.89198924712221
$ws.Range("F16").Value = 24.24618066431346
$ws.Range("G16").Value = 29.78934952597182
$ws.Range("H16").Value = 13.53938743163063
$ws.Range("I16").Value = 21.94513449073567
$ws.Range("L16").Value = 9.015711139761914
$ws.Range("M16").Value = 23.92465300765524
$ws.Range("N16").Value = 18.59134044427574
$ws.Range("O16").Value = 20.95511424748181
$ws.Range("C17").Value = 9.936599780215378
$ws.Range("D17").Value = 4.035115041386867
$ws.Range("E17").Value = 11.92312255454585
$ws.Range("F17").Value = 24.13692801119759
$ws.Range("G17").Value = 29.57005003410401
$ws.Range("H17").Value = 13.53150889176508
$ws.Range("I17").Value = 21.90723035993659
$ws.Range("L17").Value = 9.033980058976358
$ws.Range("M17").Value = 23.63653549479731
$ws.Range("N17").Value = 18.45390984390552
$ws.Range("O17").Value = 20.90218826695442
$ws.Range("C18").Value = 9.945972454487048
$ws.Range("D18").Value = 4.031143759116227
$ws.Range("E18").Value = 11.94126811003573
$ws.Range("F18").Value = 24.07457285859073
$ws.Range("G18").Value = 29.44415440853185
$ws.Range("H18").Value = 13.52726514894527
$ws.Range("I18").Value = 21.88599472864551
$ws.Range("L18").Value = 9.044617139164741
$ws.Range("M18").Value = 23.46912446418831
$ws.Range("N18").Value = 18.37454206503231
$ws.Range("O18").Value = 20.87225241321695
$ws.Range("C19").Value = 9.949182164236563
$ws.Range("D19").Value = 4.029795222421265
$ws.Range("E19").Value = 11.94745286901321
$ws.Range("F19").Value = 24.0535456232913
$ws.Range("G19").Value = 29.40157443987786
$ws.Range("H19").Value = 13.52587778861027
$ws.Range("I19").Value = 21.87890232960881
$ws.Range("L19").Value = 9.048240902290308
$ws.Range("M19").Value = 23.41215554599557
$ws.Range("N19").Value = 18.34761684564466
$ws.Range("O19").Value = 20.86220424379139
$ws.Range("C20").Value = 9.934882352336331
$ws.Range("D20").Value = 4.035848163766505
$ws.Range("E20").Value = 11.91978368233376
$ws.Range("F20").Value = 24.14850851959498
$ws.Range("G20").Value = 29.59337121566242
$ws.Range("H20").Value = 13.53231780821941
$ws.Range("I20").Value = 21.91120686463983
$ws.Range("L20").Value = 9.032021927491
$ws.Range("M20").Value = 23.66738240416817
$ws.Range("N20").Value = 18.4685734224971
$ws.Range("O20").Value = 20.90777015717303
$ws.Range("C21").Value = 9.889366549354966
$ws.Range("D21").Value = 4.055907912921009
$ws.Range("E21").Value = 11.82966690106572
$ws.Range("F21").Value = 24.47320499700623
$ws.Range("G21").Value = 30.24028172481463
$ws.Range("H21").Value = 13.55738814050736
$ws.Range("I21").Value = 22.02647196417683
$ws.Range("L21").Value = 8.979071075176574
$ws.Range("M21").Value = 24.50436453183184
$ws.Range("N21").Value = 18.8709481956658
$ws.Range("O21").Value = 21.06685273903086
$ws.Range("C22").Value = 9.861511283113357
$ws.Range("D22").Value = 4.068806745653554
$ws.Range("E22").Value = 11.77288922166789
$ws.Range("F22").Value = 24.68926885168948
$ws.Range("G22").Value = 30.66411777749119
$ws.Range("H22").Value = 13.57631775205451
$ws.Range("I22").Value = 22.10673687856407
$ws.Range("L22").Value = 8.945609987625659
$ws.Range("M22").Value = 25.03538761783586
$ws.Range("N22").Value = 19.1305149222423
$ws.Range("O22").Value = 21.17515524579794
$ws.Range("C23").Value = 9.876204624654253
$ws.Range("D23").Value = 4.061942037163786
$ws.Range("E23").Value = 11.80299889683798
$ws.Range("F23").Value = 24.57362086631534
$ws.Range("G23").Value = 30.43786495906984
$ws.Range("H23").Value = 13.56598277813892
$ws.Range("I23").Value = 22.06345208105256
$ws.Range("L23").Value = 8.963364282893901
$ws.Range("M23").Value = 24.75347510656846
$ws.Range("N23").Value = 18.99231677091065
$ws.Range("O23").Value = 21.11696426502416
$ws.Range("C24").Value = 9.935658129919885
$ws.Range("D24").Value = 4.03551679637627
$ws.Range("E24").Value = 11.92129241822947
$ws.Range("F24").Value = 24.1432715473239
$ws.Range("G24").Value = 29.58282712840773
$ws.Range("H24").Value = 13.53195120681662
$ws.Range("I24").Value = 21.9094073539054
$ws.Range("L24").Value = 9.032906781139008
$ws.Range("M24").Value = 23.65344201413283
$ws.Range("N24").Value = 18.46194512468447
$ws.Range("O24").Value = 20.90524505014969
$ws.Range("C25").Value = 10.0078687201673
$ws.Range("D25").Value = 4.006135780756886
$ws.Range("E25").Value = 12.05801307495498
$ws.Range("F25").Value = 23.69783274301685
$ws.Range("G25").Value = 28.67037421640079
$ws.Range("H25").Value = 13.5062513132575
$ws.Range("I25").Value = 21.76492463119945
$ws.Range("L25").Value = 9.112866040912149
$ws.Range("M25").Value = 22.402834813771
$ws.Range("N25").Value = 17.87793154493992
$ws.Range("O25").Value = 20.69629302094984

Write-Host "Updated loading_percent values for case with 380 kV"